$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# 1) Remove the _GoBack bookmark from the "Jump to start/end of file" cell
#    by rewriting that paragraph (preserving its original attributes) without
#    the bookmark markers.
$jumpCell = $t.Cell(2, 2).Range
$jumpXml = $pkgHeader + '<w:body><w:p w14:paraId="3063A04B" w14:textId="75434569" w:rsidR="008849DF" w:rsidRDefault="008849DF" w:rsidP="00CD01A8"><w:r><w:t>Jump to start/end of file</w:t></w:r></w:p></w:body>' + $pkgFooter
$jumpCell.InsertXML($jumpXml)

# 2) Insert a new row right after the "Navigate to class" row with the
#    Cmd+Shift+O shortcut, and move the _GoBack bookmark there.
$navClassRow = 0
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    if ($t.Cell($i, 2).Range.Text.StartsWith("Navigate to class")) {
        $navClassRow = $i
        break
    }
}

$newRow = $t.Rows.Add($t.Rows.Item($navClassRow + 1))

$shortcutCell = $t.Cell($navClassRow + 1, 1).Range
$shortcutXml = $pkgHeader + '<w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Cmd+Shift+O</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body>' + $pkgFooter
$shortcutCell.InsertXML($shortcutXml)

$descCell = $t.Cell($navClassRow + 1, 2).Range
$descXml = $pkgHeader + '<w:body><w:p><w:r><w:t>Open file dialog. You can use wildcards (*) to search with partial names.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body>' + $pkgFooter
$descCell.InsertXML($descXml)
